# Actualización automática 2025-09-01 08:30:07
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L3").Value = 0
$wsGrupo.Range("L13").Value = "0 de 11"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Column widths: C widens 10 -> 12, F widens 12 -> 16
# (ColumnWidth undergoes pixel-rounding in saved width attr; 11.1 / 15.1
#  are the character-width inputs that round-trip to exactly 12 / 16)
$wsMensual.Columns.Item(3).ColumnWidth = 11.1
$wsMensual.Columns.Item(6).ColumnWidth = 15.1

# Month headers roll forward by one month
$wsMensual.Range("C1").Value = "junio"
$wsMensual.Range("D1").Value = "julio"
$wsMensual.Range("E1").Value = "agosto"
$wsMensual.Range("F1").Value = "septiembre"

# Row 3 - ARCOS GOMEZ CONSTRUCCIONES CIA. LTDA.
$wsMensual.Range("C3").Value = 832
$wsMensual.Range("D3").Value = 594.47
$wsMensual.Range("E3").Value = 142.56
$wsMensual.Range("F3").Value = 0

# Row 5 - CARRION ALVAREZ MARIO ANDRES
$wsMensual.Range("C5").Value = 155.38
$wsMensual.Range("D5").Value = 0

# Row 11 - VACA PANCHI DORYS CAROLINA
$wsMensual.Range("C11").Value = 10.44
$wsMensual.Range("D11").Value = 0

# Row 13 - totals
$wsMensual.Range("C13").Value = 997.8200000000001
$wsMensual.Range("D13").Value = 594.47
$wsMensual.Range("E13").Value = 142.56
$wsMensual.Range("F13").Value = 0
